$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Property")

# Update IP column (F) for rows 2-6 with distinct addresses instead of the shared 127.0.0.1
$ws.Range("F2").Value = "192.168.1.113"
$ws.Range("F3").Value = "192.168.1.114"
$ws.Range("F4").Value = "192.168.1.115"
$ws.Range("F5").Value = "192.168.1.116"
$ws.Range("F6").Value = "192.168.1.117"

# Remove the extra GameServer_2 row data (row 7), keep formatting on A7:H7
$ws.Range("A7:H7").ClearContents()

# Autofit column F like Excel would do with "bestFit" (14.2857... serializes to width="15")
$ws.Columns.Item(6).ColumnWidth = 14.2857142857

# Update the selection to match the row that was just cleared
$ws.Range("A7:XFD7").Select()
